$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new header cells (I1 = "I0", J1 = "IF") ---
# Copy formatting from the existing header cell H1 (bold, bordered, centered)
# onto the new header cells, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Fill in the I0/IF data values for rows 2-67 ---
$data = @(
    @(2,6,6), @(3,9,9), @(4,7,7), @(5,6,6), @(6,7,8),
    @(7,8,8), @(8,7,7), @(9,5,5), @(10,9,9), @(11,7,8),
    @(12,6,6), @(13,7,7), @(14,7,8), @(15,8,8), @(16,5,6),
    @(17,6,6), @(18,10,10), @(19,8,8), @(20,8,8), @(21,7,7),
    @(22,8,8), @(23,6,6), @(24,5,6), @(25,8,8), @(26,5,6),
    @(27,6,6), @(28,8,8), @(29,5,5), @(30,7,8), @(31,8,8),
    @(32,8,8), @(33,7,7), @(34,6,7), @(35,7,7), @(36,9,9),
    @(37,8,8), @(38,7,7), @(39,6,6), @(40,8,8), @(41,7,7),
    @(42,9,9), @(43,8,8), @(44,6,6), @(45,6,6), @(46,6,6),
    @(47,8,8), @(48,6,7), @(49,6,6), @(50,7,7), @(51,8,8),
    @(52,5,5), @(53,6,6), @(54,5,5), @(55,9,9), @(56,7,7),
    @(57,7,7), @(58,6,6), @(59,6,6), @(60,7,7), @(61,5,5),
    @(62,6,6), @(63,7,7), @(64,6,6), @(65,7,7), @(66,9,9),
    @(67,3,3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
